# working buttons and basic screens
# - Tighten the tip-temp conversion table (rows 31:33): the ADC oversampling
#   divisor (H) drops from 100 to 10, and the final ratio (K) is now rounded
#   to the nearest integer instead of left as a raw fraction.
# - Add a new scratch area below the table (rows 38/39/42) for working out
#   the new screen's temperature-range labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: lone oversample divisor + rounded ratio ---------------------
$ws.Range("H31").Value = 10
$ws.Range("K31").Formula = "=ROUND(I31/J31, 0)"

# --- Rows 32:33: shared oversample divisor + rounded ratio ---------------
$ws.Range("H32").Value = 10
$ws.Range("H33").Value = 10
$ws.Range("K32:K33").Formula = "=ROUND(I32/J32, 0)"

# --- New scratch rows for the basic-screens temperature labels -----------
$ws.Range("B38").Value = "80C"
$ws.Range("C38").Value = "190C"
$ws.Range("B39").Formula = "=80-25"
$ws.Range("B42").Value = "500-25"

# Leave the selection where the author left off editing
[void]$ws.Range("L33").Select()
